$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width change: column 43 (AQ) from 12.7109375 to 11.7109375
# (engine quantizes ColumnWidth to 1/6-character steps; 10.877604166666666 is the
# input that lands on the closest achievable stored width to the 11.7109375 target)
$ws.Columns.Item(43).ColumnWidth = 10.877604166666666

# Cell value updates (connectivity matrix data refresh)
$ws.Range("R1").Value = 0.87539436767393952
$ws.Range("AW2").Value = 0.83133344263751441
$ws.Range("BB2").Value = 0.74519260511377827
$ws.Range("AW3").Value = 0.98215268147979484
$ws.Range("G5").Value = 0.92521522176088888
$ws.Range("U6").Value = 0.7966650663301218
$ws.Range("Z6").Value = 0.98310548653807306
$ws.Range("BI7").Value = 0.83077087692626972
$ws.Range("F8").Value = 0.71013080558450192
$ws.Range("G8").Value = 0.99554173184031314
$ws.Range("BI8").Value = 0.9950003889510064
$ws.Range("X9").Value = 0.92233831789431275
$ws.Range("AS9").Value = 0.98678195115345724
$ws.Range("C10").Value = 0.72757062326222688
$ws.Range("AI10").Value = 0.7557140317945592
$ws.Range("BE11").Value = 0.81371415214196618
$ws.Range("O12").Value = 0.89357224237951671
$ws.Range("H13").Value = 0.78970142648932806
$ws.Range("P13").Value = 0.89599349391414829
$ws.Range("BI14").Value = 0.77597898376850882
$ws.Range("T15").Value = 0.774836813223444
$ws.Range("BF15").Value = 0.99799624866911629
$ws.Range("AN16").Value = 0.82996902024675279
$ws.Range("D17").Value = 0.66094959973712686
$ws.Range("T17").Value = 0.92546093057608769
$ws.Range("Z17").Value = 0.66450740815205722
$ws.Range("BI17").Value = 0.74651903865473757
$ws.Range("T18").Value = 0.75573408384022867
$ws.Range("AC18").Value = 0.72185519272553933
$ws.Range("BB18").Value = 0.66970377954097682
$ws.Range("B19").Value = 0.9470554300551961
$ws.Range("N19").Value = 0.87301567245480094
$ws.Range("AJ21").Value = 0.85556151230552113
$ws.Range("N22").Value = 0.70339792924640809
$ws.Range("AO22").Value = 0.99767314791234929
$ws.Range("BN22").Value = 0.8964950955973382
$ws.Range("B23").Value = 0.70752144717077847
$ws.Range("L24").Value = 0.70304631439528409
$ws.Range("Z24").Value = 0.82693029387515349
$ws.Range("BD24").Value = 0.86119130717728443
$ws.Range("W26").Value = 0.99666586500842191
$ws.Range("BE26").Value = 0.86213539584725796
$ws.Range("O27").Value = 0.80878637384076391
$ws.Range("Y27").Value = 0.90949351046374405
$ws.Range("AD28").Value = 0.84415026370359647
$ws.Range("AL28").Value = 0.86014565162870971
$ws.Range("AX28").Value = 0.84565120212650846
$ws.Range("BL28").Value = 0.94861237678655541
$ws.Range("F29").Value = 0.9671931694757494
$ws.Range("Y29").Value = 0.92359509880427915
$ws.Range("AR29").Value = 0.8749302959270796
$ws.Range("C30").Value = 0.77816886796340912
$ws.Range("AF30").Value = 0.8311696256129667
$ws.Range("AT30").Value = 0.98381327150265152
$ws.Range("B31").Value = 0.60852216042830609
$ws.Range("T31").Value = 0.63882494598736717
$ws.Range("AK31").Value = 0.87568243785120414
$ws.Range("AV31").Value = 0.7983534172168284
$ws.Range("N32").Value = 0.97553825050745346
$ws.Range("N33").Value = 0.70755848477880035
$ws.Range("AH33").Value = 0.96792746317401579
$ws.Range("AV34").Value = 0.86554929560793081
$ws.Range("AR36").Value = 0.99086910591270372
$ws.Range("AB37").Value = 0.80485811521519945
$ws.Range("AI37").Value = 0.94664849488169511
$ws.Range("BE37").Value = 0.93277419453737775
$ws.Range("AP38").Value = 0.87237188530340237
$ws.Range("BK38").Value = 0.99700221058277294
$ws.Range("BB39").Value = 0.66777292563168467
$ws.Range("AD41").Value = 0.68359533035838083
$ws.Range("AN41").Value = 0.99266057540596608
$ws.Range("P42").Value = 0.92613360296132896
$ws.Range("BJ42").Value = 0.62161112822330788
$ws.Range("BM42").Value = 0.62715130580052048
$ws.Range("L43").Value = 0.6507473269922599
$ws.Range("W44").Value = 0.97090119782265749
$ws.Range("A45").Value = 0.76806242902769939
$ws.Range("AE46").Value = 0.67818748071419721
$ws.Range("AS46").Value = 0.86931768624186878
$ws.Range("AX46").Value = 0.90627475465617624
$ws.Range("BH46").Value = 0.88991949658146619
$ws.Range("BN46").Value = 0.72801482836863096
$ws.Range("AN47").Value = 0.80701519132802568
$ws.Range("Y48").Value = 0.62782083993917426
$ws.Range("AK48").Value = 0.89260452234531029
$ws.Range("E49").Value = 0.75974268652402843
$ws.Range("AJ49").Value = 0.87200383537626969
$ws.Range("AV49").Value = 0.60552367430168674
$ws.Range("BB49").Value = 0.93812135179609724
$ws.Range("BE49").Value = 0.65354163404904253
$ws.Range("A50").Value = 0.86763865991618816
$ws.Range("S50").Value = 0.59233179534163738
$ws.Range("BK51").Value = 0.9685305770408863
$ws.Range("BM51").Value = 0.76866145681978482
$ws.Range("AA52").Value = 0.86727563974467325
$ws.Range("AW52").Value = 0.81295571581994475
$ws.Range("D53").Value = 0.98341416651055358
$ws.Range("E53").Value = 0.74869632660616203
$ws.Range("N53").Value = 0.88750613676364765
$ws.Range("AY54").Value = 0.60540595364711791
$ws.Range("H55").Value = 0.68374601463297247
$ws.Range("Q55").Value = 0.91724575417825416
$ws.Range("X55").Value = 0.98515719812319036
$ws.Range("AA55").Value = 0.83746623778468998
$ws.Range("AE55").Value = 0.99942244989502216
$ws.Range("AI55").Value = 0.74160056284242259
$ws.Range("Z56").Value = 0.9044044128407962
$ws.Range("BC56").Value = 0.9633634469060538
$ws.Range("W57").Value = 0.90971213867816392
$ws.Range("AQ57").Value = 0.97218231298211033
$ws.Range("AZ57").Value = 0.82353522875643448
$ws.Range("BM57").Value = 0.87746153660944504
$ws.Range("BA58").Value = 0.9941111373831949
$ws.Range("BH58").Value = 0.89404298133815918
$ws.Range("Z59").Value = 0.57005036682888455
$ws.Range("AH59").Value = 0.8125197476069328
$ws.Range("BB59").Value = 0.96797042838159864
$ws.Range("BE59").Value = 0.7354532271914811
$ws.Range("D60").Value = 0.97106767698198304
$ws.Range("K60").Value = 0.60328677133210651
$ws.Range("AN61").Value = 0.84891147684378121
$ws.Range("BG61").Value = 0.74064544122242126
$ws.Range("G62").Value = 0.95729857729337287
$ws.Range("N62").Value = 0.78203846559170986
$ws.Range("BH62").Value = 0.77788999551573879
$ws.Range("S63").Value = 0.77924843823912804
$ws.Range("AC63").Value = 0.98765872656990039
$ws.Range("AM64").Value = 0.62968682390953135
$ws.Range("BN64").Value = 0.71939097881943737
$ws.Range("BA66").Value = 0.89945218830239049
$ws.Range("S67").Value = 0.97984073280540829
$ws.Range("AF67").Value = 0.85906897947064531
$ws.Range("AU67").Value = 0.92803161779522458
$ws.Range("BA67").Value = 0.71908858968846678
$ws.Range("B68").Value = 0.92943892555646146
$ws.Range("BI68").Value = 0.98989325816058993
